$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to stay text (matches source values like "12.99"),
# otherwise Excel auto-converts numeric-looking strings to numbers.
$ws.Range("C2:C11").NumberFormat = "@"

$data = @(
  @('Squishmallows 16" Plush', 'Squishmallow', '12.99', '16"', 'Multi', 'Yes', 'https://www.costco.com/squishmallows-16%22-plush.product.100734711.html', '8011603391.jpg', 'Costco'),
  @('Squishmallows 16" Plush', 'Squishmallow', '12.99', '16"', 'Teal', 'Yes', 'https://www.costco.com/squishmallows-16%22-plush.product.100734711.html', '7061603393.jpg', 'Costco'),
  @('Squishmallows 16" Plush', 'Squishmallow', '12.99', '16"', 'Brown', 'Yes', 'https://www.costco.com/squishmallows-16%22-plush.product.100734711.html', '361603390.jpg', 'Costco'),
  @('Squishmallows 16" Plush', 'Squishmallow', '12.99', '16"', 'Pink', 'Yes', 'https://www.costco.com/squishmallows-16%22-plush.product.100734711.html', '4401603392.jpg', 'Costco'),
  @('Squishmallows 16" Plush', 'Squishmallow', '12.99', '16"', 'Multi', 'Yes', 'https://www.costco.com/squishmallows-16%22-plush.product.100734711.html', '3541603391.jpg', 'Costco'),
  @('Squishmallows 16" Plush', 'Squishmallow', '12.99', '16"', 'Teal', 'Yes', 'https://www.costco.com/squishmallows-16%22-plush.product.100734711.html', '5221603393.jpg', 'Costco'),
  @('Squishmallows 16" Plush', 'Squishmallow', '12.99', '16"', 'Brown', 'Yes', 'https://www.costco.com/squishmallows-16%22-plush.product.100734711.html', '8231603390.jpg', 'Costco'),
  @('Squishmallows 16" Plush', 'Squishmallow', '12.99', '16"', 'Pink', 'Yes', 'https://www.costco.com/squishmallows-16%22-plush.product.100734711.html', '6141603392.jpg', 'Costco'),
  @('Squishmallows 20" Hello Kitty Sunglasses', ' Hello Kitty Sunglasses', '27.99', '20"', 'None', 'Yes', 'https://www.costco.com/squishmallows-20%22-hello-kitty-sunglasses.product.100742307.html', '1603305.jpg', 'Costco'),
  @('Squishmallows 20” Star Wars Chewbacca Plush', 'Squishmallows 20” Star Wars Chewbacca ', '27.99', 'N/A', 'None', 'Yes', 'https://www.costco.com/squishmallows-20%e2%80%9d-star-wars-chewbacca-plush.product.100691777.html', '1545490.jpg', 'Costco'),
)

$r = 2
foreach ($row in $data) {
  $c = 1
  foreach ($val in $row) {
    $ws.Cells.Item($r, $c).Value = $val
    $c = $c + 1
  }
  $r = $r + 1
}

# Restore the selection to the refreshed data block (rows 2-7 were the
# originally-selected rows in the source workbook).
$ws.Rows("2:7").Select()